$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = "primary dict"
$ws.Range("N5").Value = "crime,criminals,theft,violent crime,criminal violence,homeless,rule of law,law and order,gangs"
$ws.Range("N15").Value = "foreign policy,war,china,nato,putin,xi,jinping,sanctions,taiwan,international,world leaders,north korea,european union,israel,iran,afghanistan"
$ws.Range("N3").Value = "education,schools,school,teachers,teacher,classroom,classrooms,book bans,banned book,parents,school funding,students,age appropriate,sex ed"
$ws.Range("N11").Value = "supreme court nomination,supreme court appointee,justices,appoint judges"
